$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Cxcl10"
$ws.Cells.Item(2,3).Value = "Ccr3"
$ws.Cells.Item(2,4).Value = "M1"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 6.643695666666666
$ws.Cells.Item(2,8).Value = 19.931087
$ws.Cells.Item(2,9).Value = 0.02330062621916135
$ws.Cells.Item(2,10).Value = 0.0235018519445706
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.1323866666666667
$ws.Cells.Item(2,14).Value = 0.39716
$ws.Cells.Item(2,15).Value = 0.2977240407890032
$ws.Cells.Item(2,16).Value = 0.2977240407890032
$ws.Cells.Item(2,17).Value = 0.8795367236577776
$ws.Cells.Item(2,18).Value = 7.915830512919999
$ws.Cells.Item(2,19).Value = 0.00693715659088291
$ws.Cells.Item(2,20).Value = 0.006997066326962451

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Cxcl10"
$ws.Cells.Item(3,3).Value = "Ccr3"
$ws.Cells.Item(3,4).Value = "M2"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 6.643695666666666
$ws.Cells.Item(3,8).Value = 19.931087
$ws.Cells.Item(3,9).Value = 0.02330062621916135
$ws.Cells.Item(3,10).Value = 0.0235018519445706
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.3122756666666667
$ws.Cells.Item(3,14).Value = 0.936827
$ws.Cells.Item(3,15).Value = 0.7022759592109968
$ws.Cells.Item(3,16).Value = 0.7022759592109967
$ws.Cells.Item(3,17).Value = 2.074664493438778
$ws.Cells.Item(3,18).Value = 18.671980440949
$ws.Cells.Item(3,19).Value = 0.01636346962827844
$ws.Cells.Item(3,20).Value = 0.01650478561760815

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Cxcl10"
$ws.Cells.Item(4,3).Value = "Ccr3"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 139.9983773333333
$ws.Cells.Item(4,8).Value = 419.995132
$ws.Cells.Item(4,9).Value = 0.490999290936783
$ws.Cells.Item(4,10).Value = 0.4952395927881098
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.1323866666666667
$ws.Cells.Item(4,14).Value = 0.39716
$ws.Cells.Item(4,15).Value = 0.2977240407890032
$ws.Cells.Item(4,16).Value = 0.2977240407890032
$ws.Cells.Item(4,17).Value = 18.53391851390222
$ws.Cells.Item(4,18).Value = 166.80526662512
$ws.Cells.Item(4,19).Value = 0.1461822929222344
$ws.Cells.Item(4,20).Value = 0.1474447327235765

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Cxcl10"
$ws.Cells.Item(5,3).Value = "Ccr3"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 139.9983773333333
$ws.Cells.Item(5,8).Value = 419.995132
$ws.Cells.Item(5,9).Value = 0.490999290936783
$ws.Cells.Item(5,10).Value = 0.4952395927881098
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.3122756666666667
$ws.Cells.Item(5,14).Value = 0.936827
$ws.Cells.Item(5,15).Value = 0.7022759592109968
$ws.Cells.Item(5,16).Value = 0.7022759592109967
$ws.Cells.Item(5,17).Value = 43.71808661401823
$ws.Cells.Item(5,18).Value = 393.462779526164
$ws.Cells.Item(5,19).Value = 0.3448169980145486
$ws.Cells.Item(5,20).Value = 0.3477948600645332

# Row 6
$ws.Cells.Item(6,1).Value = "M1"
$ws.Cells.Item(6,2).Value = "Cxcl10"
$ws.Cells.Item(6,3).Value = "Ccr3"
$ws.Cells.Item(6,4).Value = "M1"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 92.13145400000001
$ws.Cells.Item(6,8).Value = 276.394362
$ws.Cells.Item(6,9).Value = 0.3231214493241425
$ws.Cells.Item(6,10).Value = 0.3259119472027818
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.1323866666666667
$ws.Cells.Item(6,14).Value = 0.39716
$ws.Cells.Item(6,15).Value = 0.2977240407890032
$ws.Cells.Item(6,16).Value = 0.2977240407890032
$ws.Cells.Item(6,17).Value = 12.19697609021333
$ws.Cells.Item(6,18).Value = 109.77278481192
$ws.Cells.Item(6,19).Value = 0.09620102355838282
$ws.Cells.Item(6,20).Value = 0.09703182186262446

# Row 7
$ws.Cells.Item(7,1).Value = "M1"
$ws.Cells.Item(7,2).Value = "Cxcl10"
$ws.Cells.Item(7,3).Value = "Ccr3"
$ws.Cells.Item(7,4).Value = "M2"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 92.13145400000001
$ws.Cells.Item(7,8).Value = 276.394362
$ws.Cells.Item(7,9).Value = 0.3231214493241425
$ws.Cells.Item(7,10).Value = 0.3259119472027818
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.3122756666666667
$ws.Cells.Item(7,14).Value = 0.936827
$ws.Cells.Item(7,15).Value = 0.7022759592109968
$ws.Cells.Item(7,16).Value = 0.7022759592109967
$ws.Cells.Item(7,17).Value = 28.77041121881934
$ws.Cells.Item(7,18).Value = 258.933700969374
$ws.Cells.Item(7,19).Value = 0.2269204257657597
$ws.Cells.Item(7,20).Value = 0.2288801253401573

# Row 8
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Cxcl10"
$ws.Cells.Item(8,3).Value = "Ccr3"
$ws.Cells.Item(8,4).Value = "M1"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 39.03202233333334
$ws.Cells.Item(8,8).Value = 117.096067
$ws.Cells.Item(8,9).Value = 0.1368922672858171
$ws.Cells.Item(8,10).Value = 0.1380744778207792
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.1323866666666667
$ws.Cells.Item(8,14).Value = 0.39716
$ws.Cells.Item(8,15).Value = 0.2977240407890032
$ws.Cells.Item(8,16).Value = 0.2977240407890032
$ws.Cells.Item(8,17).Value = 5.167319329968889
$ws.Cells.Item(8,18).Value = 46.50587396972
$ws.Cells.Item(8,19).Value = 0.04075611896910174
$ws.Cells.Item(8,20).Value = 0.04110809146663397

# Row 9
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Cxcl10"
$ws.Cells.Item(9,3).Value = "Ccr3"
$ws.Cells.Item(9,4).Value = "M2"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 39.03202233333334
$ws.Cells.Item(9,8).Value = 117.096067
$ws.Cells.Item(9,9).Value = 0.1368922672858171
$ws.Cells.Item(9,10).Value = 0.1380744778207792
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.3122756666666667
$ws.Cells.Item(9,14).Value = 0.936827
$ws.Cells.Item(9,15).Value = 0.7022759592109968
$ws.Cells.Item(9,16).Value = 0.7022759592109967
$ws.Cells.Item(9,17).Value = 12.18875079548989
$ws.Cells.Item(9,18).Value = 109.698757159409
$ws.Cells.Item(9,19).Value = 0.09613614831671537
$ws.Cells.Item(9,20).Value = 0.09696638635414519

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Cxcl10"
$ws.Cells.Item(10,3).Value = "Ccr3"
$ws.Cells.Item(10,4).Value = "M1"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 7.3239405
$ws.Cells.Item(10,8).Value = 14.647881
$ws.Cells.Item(10,9).Value = 0.02568636623409617
$ws.Cells.Item(10,10).Value = 0.01727213024375885
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.1323866666666667
$ws.Cells.Item(10,14).Value = 0.39716
$ws.Cells.Item(10,15).Value = 0.2977240407890032
$ws.Cells.Item(10,16).Value = 0.2977240407890032
$ws.Cells.Item(10,17).Value = 0.9695920696599999
$ws.Cells.Item(10,18).Value = 5.817552417959999
$ws.Cells.Item(10,19).Value = 0.007647448748401322
$ws.Cells.Item(10,20).Value = 0.005142328409205834

# Row 11
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Cxcl10"
$ws.Cells.Item(11,3).Value = "Ccr3"
$ws.Cells.Item(11,4).Value = "M2"
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 7.3239405
$ws.Cells.Item(11,8).Value = 14.647881
$ws.Cells.Item(11,9).Value = 0.02568636623409617
$ws.Cells.Item(11,10).Value = 0.01727213024375885
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.3122756666666667
$ws.Cells.Item(11,14).Value = 0.936827
$ws.Cells.Item(11,15).Value = 0.7022759592109968
$ws.Cells.Item(11,16).Value = 0.7022759592109967
$ws.Cells.Item(11,17).Value = 2.2870884022645
$ws.Cells.Item(11,18).Value = 13.722530413587
$ws.Cells.Item(11,19).Value = 0.01803891748569485
$ws.Cells.Item(11,20).Value = 0.01212980183455301
